$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (shared string referenced by A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 15:11"

# Per-cell updates derived from the country data refresh + re-sort (see commit diff)
# Row 4
$ws.Range("B4").Value = 3292681
$ws.Range("C4").Value = 895
$ws.Range("D4").Value = 1460649
$ws.Range("E4").Value = 1695316
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = 136716
# Row 6
$ws.Range("B6").Value = 825736
$ws.Range("C6").Value = 3133
$ws.Range("D6").Value = 517546
$ws.Range("E6").Value = 286019
$ws.Range("G6").Value = 27
$ws.Range("H6").Value = 22171
# Row 17
$ws.Range("B17").Value = 229480
$ws.Range("C17").Value = 2994
$ws.Range("D17").Value = 165396
$ws.Range("E17").Value = 61903
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 2181
# Row 25
$ws.Range("D25").Value = 41408
$ws.Range("E25").Value = 50865
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 1787
# Row 40
$ws.Range("B40").Value = 50921
$ws.Range("C40").Value = 81
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 6137
# Row 58
$ws.Range("B58").Value = 23521
$ws.Range("C58").Value = 531
$ws.Range("D58").Value = 14607
$ws.Range("E58").Value = 8616
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 298
# Row 63
$ws.Range("B63").Value = 18073
$ws.Range("C63").Value = 345
$ws.Range("D63").Value = 13780
$ws.Range("E63").Value = 3911
$ws.Range("G63").Value = 12
$ws.Range("H63").Value = 382
# Row 70
$ws.Range("B70").Value = 12402
$ws.Range("C70").Value = 375
$ws.Range("E70").Value = 4806
# Row 87
$ws.Range("A87").Value = "Bosnia y Herzegovina"
$ws.Range("B87").Value = 6719
$ws.Range("C87").Value = 317
$ws.Range("D87").Value = 3078
$ws.Range("E87").Value = 3422
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 219
# Row 88
$ws.Range("A88").Value = "Haiti"
$ws.Range("B88").Value = 6617
$ws.Range("C88").Value = 35
$ws.Range("D88").Value = 2590
$ws.Range("E88").Value = 3892
$ws.Range("G88").Value = 5
$ws.Range("H88").Value = 135
# Row 89
$ws.Range("A89").Value = "Tayikistan"
$ws.Range("B89").Value = 6457
$ws.Range("D89").Value = 5115
$ws.Range("E89").Value = 1287
$ws.Range("H89").Value = 55
# Row 101
$ws.Range("B101").Value = 3672
$ws.Range("C101").Value = 140
$ws.Range("D101").Value = 2466
$ws.Range("E101").Value = 1088
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 118
# Row 142
$ws.Range("A142").Value = "Liberia"
$ws.Range("B142").Value = 998
$ws.Range("C142").Value = 35
$ws.Range("D142").Value = 420
$ws.Range("E142").Value = 531
$ws.Range("H142").Value = 47
# Row 143
$ws.Range("A143").Value = "Uruguay"
$ws.Range("B143").Value = 985
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 886
$ws.Range("E143").Value = 70
$ws.Range("H143").Value = 29
# Row 144
$ws.Range("A144").Value = "Georgia"
$ws.Range("B144").Value = 981
$ws.Range("C144").Value = 8
$ws.Range("D144").Value = 851
$ws.Range("E144").Value = 115
$ws.Range("H144").Value = 15
# Row 209
$ws.Range("A209").Value = "Groenlandia"
# Row 210
$ws.Range("A210").Value = "Islas Malvinas"
# Row 211
$ws.Range("A211").Value = "Montserrat"
$ws.Range("C211").Value = 1
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 1
$ws.Range("H211").Value = 1
# Row 212
$ws.Range("A212").Value = "Santa Sede"
$ws.Range("B212").Value = 12
$ws.Range("D212").Value = 12
$ws.Range("E212").Value = 0
# Row 213
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("E213").Value = 3
$ws.Range("H213").Value = 0
